# backbone v15.3 in all functions
# Fill column V (rows 2-22) with "x" to match the rest of the row pattern,
# and move the active selection to V8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Feuil1")

# Copy the formatting used throughout column U (rows 2-22) into column V
# so the new cells match the existing row styling (centered, wrapped text).
$ws.Range("U2:U22").Copy() | Out-Null
$ws.Range("V2:V22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 22).Value = "x"
}

$ws.Range("V8").Select()
